$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder rows 3, 5, 6, 7 as part of grouping in dfg detection:
# Row3: stream:datastream/dict -> time:timestamp/datetime
# Row5: time:timestamp/datetime -> SubProcessID/str
# Row6: concept:name/str -> stream:datastream/dict
# Row7: SubProcessID/str -> concept:name/str

$ws.Range("A3").Value = "time:timestamp"
$ws.Range("B3").Value = "datetime"

$ws.Range("A5").Value = "SubProcessID"
$ws.Range("B5").Value = "str"

$ws.Range("A6").Value = "stream:datastream"
$ws.Range("B6").Value = "dict"

$ws.Range("A7").Value = "concept:name"
$ws.Range("B7").Value = "str"
